$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.097.20'
$ws.Range("E2").Value = '  -0.33%  '

$ws.Range("D3").Value = '2.593.42'
$ws.Range("E3").Value = '  -0.77%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").Value = "'522.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.43%  '

$ws.Range("D6").Value = "'143.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.40%  '

$ws.Range("D8").Value = "'0.567"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("D9").Value = '2.613.02'
$ws.Range("E9").Value = '  -0.40%  '

$ws.Range("D10").Value = "'6.51"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.07%  '

$ws.Range("E11").Value = '  -1.72%  '

$ws.Range("D12").Value = "'0.337"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.26%  '

$ws.Range("E13").Value = '  -0.22%  '

$ws.Range("D14").Value = '3.052.08'
$ws.Range("E14").Value = '  -0.64%  '

$ws.Range("D15").Value = '58.090.04'
$ws.Range("E15").Value = '  -0.28%  '

$ws.Range("D16").Value = "'20.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.62%  '

$ws.Range("E17").Value = '  -1.67%  '

$ws.Range("D18").Value = '2.548.11'
$ws.Range("E18").Value = '  -2.24%  '

$ws.Range("D19").Value = "'339.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.90%  '

$ws.Range("D20").Value = "'4.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.16%  '

$ws.Range("D21").Value = "'10.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.10%  '

$ws.Range("D22").Value = "'6.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.99%  '

$ws.Range("E23").Value = '  +0.10%  '

$ws.Range("D24").Value = "'65.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.09%  '

$ws.Range("E26").Value = '  -2.95%  '

$ws.Range("D27").Value = '2.715.83'
$ws.Range("E27").Value = '  -0.49%  '

$ws.Range("D29").Value = "'7.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.74%  '

$ws.Range("D30").Value = '0.0₃0746'
$ws.Range("E30").Value = '  -5.58%  '

$ws.Range("E31").Value = '  -0.03%  '

$ws.Range("D32").Value = "'6.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.87%  '

$ws.Range("D33").Value = "'1.58"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.32%  '

$ws.Range("D34").Value = "'18.76"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.04%  '

$ws.Range("D35").Value = "'149.76"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.11%  '

$ws.Range("D36").Value = "'4.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.34%  '

$ws.Range("D37").Value = "'1.13"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.36%  '

$ws.Range("D38").Value = "'0.870"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.87%  '

$ws.Range("D39").Value = "'0.863"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.76%  '

$ws.Range("D40").Value = "'35.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.94%  '

$ws.Range("D41").Value = "'1.45"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.36%  '

$ws.Range("D42").Value = "'3.54"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.69%  '

$ws.Range("D43").Value = "'0.997"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.28%  '

$ws.Range("D44").Value = "'272.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.99%  '

$ws.Range("D45").Value = "'0.600"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.11%  '

$ws.Range("D46").Value = "'0.0958"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.94%  '

$ws.Range("E47").Value = '  +0.58%  '

$ws.Range("D48").Value = "'18.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.85%  '

$ws.Range("D49").Value = "'0.0522"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.97%  '

$ws.Range("D50").Value = "'18.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.94%  '

$ws.Range("D51").Value = '1.975.21'
$ws.Range("E51").Value = '  -2.88%  '
